# TASK_6, add users table
#
# - Task "1.5 - Роутер (ApiRouter)" (row 6 on the "Задачи" sheet) gets its
#   completion ("заврешенно") date filled in (column E).
# - Task "2.1 - Таблица пользователей. (DB.users)" (row 7, the "users table"
#   task) gets its start ("начало") date filled in (column D).
# - The "Задачи" sheet becomes the active/selected sheet & tab (it was the
#   "Бэклог задач" sheet before), with cell D7 selected.
# - The "Бэклог задач" sheet is no longer the selected tab; its last
#   selection moves from C15 to C14.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Задачи")
$ws2 = $wb.Worksheets.Item("Бэклог задач")

# --- Fill in the two new date cells, reusing the existing date/time
# number format (style index 2) already used throughout these columns by
# copying formatting from a neighboring cell that already has it. ---
$ws1.Range("D6").Copy()
$ws1.Range("E6").PasteSpecial(-4122)  # xlPasteFormats
$ws1.Range("E6").Value = 41975.571527777778

$ws1.Range("C7").Copy()
$ws1.Range("D7").PasteSpecial(-4122)  # xlPasteFormats
$ws1.Range("D7").Value = 41975.571527777778

$excel.CutCopyMode = 0

# --- Update selections on both sheets, and switch the active tab from
# "Бэклог задач" to "Задачи". ---
$ws2.Range("C14").Select()

$ws1.Activate()
$ws1.Range("D7").Select()
